$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Add new row 7 to the "Heap Report from Test" sheet ---

# A7: date/time value, formatted like A2:A6 (m/d/yyyy h:mm)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 43412.73333333333

# B7, C7, D7: text values matching existing entries used in rows 4-6
$ws.Range("B7").Value = "Laptop"
$ws.Range("C7").Value = "Release"
$ws.Range("D7").Value = "Factory_Class"

# F7, G7, H7: numeric values formatted like F6:H6 (#,##0)
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = 16

$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = 16

$ws.Range("H6").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = 16

# I7: new observation text (becomes a new shared string entry)
$ws.Range("I7").Value = "Consolidated creation of static Push code snippets"

# --- Reset the view: scroll/selection back to A1 (was topLeftCell B1 / I10 selection) ---
$ws.Activate()
$ws.Range("A1").Select()
